# Add a "Ratio staff / mentor" column (F) and an "All Engineers / All Staffs"
# breakdown (rows 11-13) to the Pewlett-Hackard mentorship summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (A1:E1): switch alignment from plain "center" to
#    "centerContinuous" + vertically centered, and make the row taller so the
#    new wrapped F1 header fits.
# ---------------------------------------------------------------------------
$ws.Range("A1:E1").HorizontalAlignment = 7       # xlCenterContinuous
$ws.Range("A1:F1").VerticalAlignment = -4108     # xlCenter
$ws.Rows(1).RowHeight = 34

$ws.Columns("F").ColumnWidth = 11

# ---------------------------------------------------------------------------
# 2. New "Ratio staff / mentor" formulas down column F (rows 2-9), formatted
#    as whole numbers. F8 divides by zero (no mentors that year) so it is
#    hard-set to 0 instead of keeping the #DIV/0! formula result.
# ---------------------------------------------------------------------------
$ws.Range("F2").Formula = "=(D2-C2-B2)/C2"
$ws.Range("F3:F9").Formula = "=(D3-C3-B3)/C3"
$ws.Range("F8").Value = 0
$ws.Range("F2:F8").NumberFormat = "0"

# F9 (totals row) also needs the top border used by the rest of row 9.
$ws.Range("A9").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Formula = "=(D9-C9-B9)/C9"
$ws.Range("F9").NumberFormat = "0"

# Row 10 stays empty except for carrying the same number format down.
$ws.Range("F9").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").ClearContents()

# ---------------------------------------------------------------------------
# 3. Rows 11-13: engineers vs. staff roll-up with the same five metrics.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "All Engineers"
$ws.Range("B11").Formula = "=B2+B4+B7"
$ws.Range("C11:D11").Formula = "=C2+C4+C7"
$ws.Range("D11").Formula = "=D2+D4+D7"
$ws.Range("E11:E13").Formula = "=B11/D11"
$ws.Range("F11:F13").Formula = "=(D11-C11-B11)/C11"

$ws.Range("A12").Value = "All Staffs"
$ws.Range("B12").Formula = "=B3+B5+B6+B8"
$ws.Range("C12:D12").Formula = "=C3+C5+C6+C8"
$ws.Range("D12").Formula = "=D3+D5+D6+D8"

$ws.Range("B13").Formula = "=SUM(B11:B12)"
$ws.Range("C13:D13").Formula = "=SUM(C11:C12)"
$ws.Range("D13").Formula = "=SUM(D11:D12)"

# New F1 header cell - start from the existing header formatting (bold +
# underline font) then switch to a plain centered, wrapped look. Written
# after the "All Engineers" / "All Staffs" labels above so the shared-string
# table ends up in the same order as the rest of the new text.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)              # xlPasteFormats
$ws.Range("F1").Value = "Ratio staff / mentor"
$ws.Range("F1").HorizontalAlignment = -4108      # xlCenter
$ws.Range("F1").VerticalAlignment = -4108        # xlCenter
$ws.Range("F1").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Formatting for rows 11-13 - copied from the matching rows 2 / 9 so the
#    new rows look identical to the existing table (comma / percent / plain
#    integer number formats, and the bordered "totals" look on row 13).
# ---------------------------------------------------------------------------
$ws.Range("B2:D2").Copy()
$ws.Range("B11:D11").PasteSpecial(-4122)
$ws.Range("B12:D12").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E11:E12").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("F11:F12").PasteSpecial(-4122)

# Row 13 (totals) reuses row 9's bordered formatting.
$ws.Range("A9").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("B9:D9").Copy()
$ws.Range("B13:D13").PasteSpecial(-4122)

$ws.Range("E9").Copy()
$ws.Range("E13").PasteSpecial(-4122)

$ws.Range("F9").Copy()
$ws.Range("F13").PasteSpecial(-4122)

# Re-apply the whole-number format to the F column after the formatting
# copy/paste passes above (PasteSpecial can drag along the source's percent /
# comma format on adjoining cells).
$ws.Range("F2:F13").NumberFormat = "0"

$ws.Range("A1").Select()
